# Update odds figures for the week of 2025-04-22 FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Manchester City vs Aston Villa
$ws.Range("J2").Value = 1.02
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 1.14
$ws.Range("M2").Value = 5.5

# Row 4 - Valencia vs Espanyol
$ws.Range("L4").Value = 1.36
$ws.Range("M4").Value = 3.2

# Row 5 - Barcelona vs Mallorca
$ws.Range("AH5").Value = 126

# Row 17 - HJK vs Gnistan
$ws.Range("G17").Value = 1.45
$ws.Range("N17").Value = 1.62

# Row 18 - Derry City vs Sligo Rovers
$ws.Range("G18").Value = 1.48
$ws.Range("J18").Value = 1.05
$ws.Range("L18").Value = 1.25

# Row 19 - Stal Mielec vs Gornik Zabrze
$ws.Range("J19").Value = 1.05
$ws.Range("L19").Value = 1.25
$ws.Range("R19").Value = 1.67

# Row 20 - Kotwica Kolobrzeg vs Ruch Chorzow
$ws.Range("J20").Value = 1.06
$ws.Range("L20").Value = 1.3
$ws.Range("R20").Value = 1.87
$ws.Range("S20").Value = 1.77

# Row 21 - Damac vs Al Nassr
$ws.Range("L21").Value = 1.17

# Row 22 - Al Wehda vs Al Ahli SC
$ws.Range("J22").Value = 1.02
$ws.Range("L22").Value = 1.13
